# Auto-generated edit script applying numeric corrections to the Lamia_Profits
# leve-profit data across all 8 job sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
# For each affected row, columns H-N (currentAveragePrice*, LevePrice*, LeveProfit*)
# are updated to their new values. A couple of cells are removed or newly added
# where the corresponding <c> element appeared/disappeared in the source diff.

$wb = $excel.ActiveWorkbook

# ===== Sheet: ALC =====
$ws = $wb.Worksheets.Item("ALC")

# Row 32
$ws.Range("H32").Value = 9498.625
$ws.Range("I32").Value = 7447
$ws.Range("J32").Value = 10182.5
$ws.Range("K32").Value = 7447
$ws.Range("L32").Value = 10182.5
$ws.Range("M32").Value = -7121
$ws.Range("N32").Value = -10834.5

# Row 112
$ws.Range("H112").Value = 1631.3334
$ws.Range("J112").Value = 1703
$ws.Range("L112").Value = 5109
$ws.Range("N112").Value = -7325

# Row 137
$ws.Range("H137").Value = 3090.8096
$ws.Range("J137").Value = 4039.818
$ws.Range("L137").Value = 12119.454
$ws.Range("N137").Value = -17219.454

# Row 138
$ws.Range("H138").Value = 3776.575
$ws.Range("J138").Value = 3818.0715
$ws.Range("L138").Value = 11454.2145
$ws.Range("N138").Value = -21734.2145

# Row 141
$ws.Range("H141").Value = 5269.676
$ws.Range("I141").Value = 3639.12
$ws.Range("K141").Value = 10917.36
$ws.Range("M141").Value = -5737.360000000001

# ===== Sheet: ARM =====
$ws = $wb.Worksheets.Item("ARM")

# Row 32
$ws.Range("H32").Value = 2269.5894
$ws.Range("I32").Value = 1794.44
$ws.Range("J32").Value = 6229.1665
$ws.Range("K32").Value = 1794.44
$ws.Range("L32").Value = 6229.1665
$ws.Range("M32").Value = -1507.44
$ws.Range("N32").Value = -6803.1665

# Row 61
$ws.Range("H61").Value = 4310.9834
$ws.Range("I61").Value = 3271.6316
$ws.Range("J61").Value = 6028.174
$ws.Range("K61").Value = 3271.6316
$ws.Range("L61").Value = 6028.174
$ws.Range("M61").Value = -3059.6316
$ws.Range("N61").Value = -6452.174

# Row 74
$ws.Range("H74").Value = 5294942
$ws.Range("I74").Value = 6948319
$ws.Range("J74").Value = 4136.6665
$ws.Range("K74").Value = 6948319
$ws.Range("L74").Value = 4136.6665
$ws.Range("M74").Value = -6947445
$ws.Range("N74").Value = -5884.6665

# Row 77
$ws.Range("H77").Value = 5294942
$ws.Range("I77").Value = 6948319
$ws.Range("J77").Value = 4136.6665
$ws.Range("K77").Value = 34741595
$ws.Range("L77").Value = 20683.3325
$ws.Range("M77").Value = -34737227
$ws.Range("N77").Value = -29419.3325

# Row 112
$ws.Range("H112").Value = 74375
$ws.Range("J112").Value = 74375
$ws.Range("L112").Value = 74375
$ws.Range("N112").Value = -77329

# Row 132
$ws.Range("H132").Value = 4143.909
$ws.Range("I132").Value = 2502.4707
$ws.Range("K132").Value = 7507.4121
$ws.Range("M132").Value = -4977.4121

# Row 136
$ws.Range("H136").Value = 4310.9834
$ws.Range("I136").Value = 3271.6316
$ws.Range("J136").Value = 6028.174
$ws.Range("K136").Value = 9814.8948
$ws.Range("L136").Value = 18084.522
$ws.Range("M136").Value = -7264.8948
$ws.Range("N136").Value = -23184.522

# ===== Sheet: BSM =====
$ws = $wb.Worksheets.Item("BSM")

# Row 7
$ws.Range("H7").Value = 2050
$ws.Range("I7").Value = 0
$ws.Range("K7").Value = 0
$ws.Range("M7").ClearContents()

# Row 10
$ws.Range("H10").Value = 1912.5714
$ws.Range("I10").Value = 1748
$ws.Range("J10").Value = 2900
$ws.Range("K10").Value = 1748
$ws.Range("L10").Value = 2900
$ws.Range("M10").Value = -1608
$ws.Range("N10").Value = -3180

# Row 11
$ws.Range("H11").Value = 1707.4
$ws.Range("I11").Value = 1592.6666
$ws.Range("K11").Value = 1592.6666
$ws.Range("M11").Value = -1452.6666

# Row 24
$ws.Range("H24").Value = 4633
$ws.Range("I24").Value = 4449.5
$ws.Range("K24").Value = 4449.5
$ws.Range("M24").Value = -4214.5

# Row 29
$ws.Range("H29").Value = 10666.667
$ws.Range("I29").Value = 10666.667
$ws.Range("J29").Value = 0
$ws.Range("K29").Value = 10666.667
$ws.Range("L29").Value = 0
$ws.Range("M29").Value = -10377.667
$ws.Range("N29").ClearContents()

# Row 81
$ws.Range("H81").Value = 49037.715
$ws.Range("J81").Value = 49037.715
$ws.Range("L81").Value = 49037.715
$ws.Range("N81").Value = -51159.715

# Row 84
$ws.Range("H84").Value = 49037.715
$ws.Range("J84").Value = 49037.715
$ws.Range("L84").Value = 147113.145
$ws.Range("N84").Value = -157721.145

# Row 134
$ws.Range("H134").Value = 1959.2325
$ws.Range("I134").Value = 1208.6296
$ws.Range("K134").Value = 3625.8888
$ws.Range("M134").Value = -1090.8888

# Row 137
$ws.Range("H137").Value = 71599.2
$ws.Range("J137").Value = 71599.2
$ws.Range("L137").Value = 71599.2
$ws.Range("N137").Value = -81799.2

# ===== Sheet: CRP =====
$ws = $wb.Worksheets.Item("CRP")

# Row 31
$ws.Range("H31").Value = 46663.46
$ws.Range("I31").Value = 3859.4167
$ws.Range("J31").Value = 83352.64
$ws.Range("K31").Value = 3859.4167
$ws.Range("L31").Value = 83352.64
$ws.Range("M31").Value = -3564.4167
$ws.Range("N31").Value = -83942.64

# Row 34
$ws.Range("H34").Value = 46663.46
$ws.Range("I34").Value = 3859.4167
$ws.Range("J34").Value = 83352.64
$ws.Range("K34").Value = 3859.4167
$ws.Range("L34").Value = 83352.64
$ws.Range("M34").Value = -3657.4167
$ws.Range("N34").Value = -83756.64

# Row 132
$ws.Range("H132").Value = 4048.7954
$ws.Range("I132").Value = 3887.0952
$ws.Range("J132").Value = 4196.4346
$ws.Range("K132").Value = 11661.2856
$ws.Range("L132").Value = 12589.3038
$ws.Range("M132").Value = -9131.285600000001
$ws.Range("N132").Value = -17649.3038

# Row 134
$ws.Range("H134").Value = 2971.9644
$ws.Range("I134").Value = 2483.4119
$ws.Range("J134").Value = 3727
$ws.Range("K134").Value = 7450.2357
$ws.Range("L134").Value = 11181
$ws.Range("M134").Value = -4915.2357
$ws.Range("N134").Value = -16251

# ===== Sheet: CUL =====
$ws = $wb.Worksheets.Item("CUL")

# Row 11
$ws.Range("H11").Value = 4059.8
$ws.Range("J11").Value = 8500
$ws.Range("L11").Value = 25500
$ws.Range("N11").Value = -25780

# Row 93
$ws.Range("H93").Value = 23675.666
$ws.Range("J93").Value = 23675.666
$ws.Range("L93").Value = 71026.99800000001
$ws.Range("N93").Value = -74770.99800000001

# Row 128
$ws.Range("H128").Value = 201985.67
$ws.Range("I128").Value = 201985.67
$ws.Range("K128").Value = 605957.01
$ws.Range("M128").Value = -600977.01

# Row 131
$ws.Range("H131").Value = 8103666
$ws.Range("I131").Value = 2262.375
$ws.Range("K131").Value = 6787.125
$ws.Range("M131").Value = -1747.125

# ===== Sheet: GSM =====
$ws = $wb.Worksheets.Item("GSM")

# Row 132
$ws.Range("H132").Value = 22336.076
$ws.Range("I132").Value = 32501.576
$ws.Range("J132").Value = 4680.2104
$ws.Range("K132").Value = 97504.728
$ws.Range("L132").Value = 14040.6312
$ws.Range("M132").Value = -94974.728
$ws.Range("N132").Value = -19100.6312

# Row 134
$ws.Range("H134").Value = 61499.5
$ws.Range("J134").Value = 61499.5
$ws.Range("L134").Value = 184498.5
$ws.Range("N134").Value = -189568.5

# ===== Sheet: LTW =====
$ws = $wb.Worksheets.Item("LTW")

# Row 11
$ws.Range("H11").Value = 1000
$ws.Range("J11").Value = 1000
$ws.Range("L11").Value = 1000
$ws.Range("N11").Value = -1280

# Row 19
$ws.Range("H19").Value = 5301.5
$ws.Range("I19").Value = 603
$ws.Range("J19").Value = 10000
$ws.Range("K19").Value = 603
$ws.Range("L19").Value = 10000
$ws.Range("M19").Value = -433
$ws.Range("N19").Value = -10340

# Row 23
$ws.Range("H23").Value = 30000
$ws.Range("I23").Value = 30000
$ws.Range("K23").Value = 30000
$ws.Range("M23").Value = -29770

# Row 25
$ws.Range("H25").Value = 5727.1816
$ws.Range("J25").Value = 6333.1665
$ws.Range("L25").Value = 6333.1665
$ws.Range("N25").Value = -6793.1665

# Row 104
$ws.Range("H104").Value = 31777.25
$ws.Range("J104").Value = 31777.25
$ws.Range("L104").Value = 31777.25
$ws.Range("N104").Value = -38765.25

# ===== Sheet: WVR =====
$ws = $wb.Worksheets.Item("WVR")

# Row 103
$ws.Range("H103").Value = 35000
$ws.Range("J103").Value = 35000
$ws.Range("L103").Value = 35000
$ws.Range("N103").Value = -37344

# Row 132
$ws.Range("H132").Value = 3675.2
$ws.Range("I132").Value = 3020.125
$ws.Range("J132").Value = 4839.778
$ws.Range("K132").Value = 9060.375
$ws.Range("L132").Value = 14519.334
$ws.Range("M132").Value = -6530.375
$ws.Range("N132").Value = -19579.334

# Row 137
$ws.Range("H137").Value = 70282.55499999999
$ws.Range("J137").Value = 70282.55499999999
$ws.Range("L137").Value = 70282.55499999999
$ws.Range("N137").Value = -80482.55499999999

# Row 140
$ws.Range("H140").Value = 60833
$ws.Range("J140").Value = 60833
$ws.Range("L140").Value = 60833
$ws.Range("N140").Value = -71193

